$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2028301886792453
$ws.Range("C2").Value = 0.5094339622641509
$ws.Range("J2").Value = 0.009433962264150943
$ws.Range("P2").Value = 0.1509433962264151
$ws.Range("S2").Value = 0.1273584905660377
$ws.Range("B3").Value = 0.008474576271186441
$ws.Range("C3").Value = 0.07627118644067797
$ws.Range("J3").Value = 0.02542372881355932
$ws.Range("P3").Value = 0.6779661016949152
$ws.Range("S3").Value = 0.211864406779661
$ws.Range("J4").Value = 0.03333333333333333
$ws.Range("P4").Value = 0.6666666666666666
$ws.Range("S4").Value = 0.3
$ws.Range("P5").Value = 1
$ws.Range("B6").Value = 0.02755905511811024
$ws.Range("D6").Value = 0.01181102362204724
$ws.Range("E6").Value = 0.003937007874015748
$ws.Range("F6").Value = 0.06299212598425197
$ws.Range("J6").Value = 0.2519685039370079
$ws.Range("O6").Value = 0.02362204724409449
$ws.Range("Q6").Value = 0.1692913385826772
$ws.Range("R6").Value = 0.1062992125984252
$ws.Range("S6").Value = 0.3425196850393701
$ws.Range("B7").Value = 0.1019417475728155
$ws.Range("D7").Value = 0.03883495145631068
$ws.Range("F7").Value = 0.06310679611650485
$ws.Range("J7").Value = 0.1067961165048544
$ws.Range("O7").Value = 0.03398058252427184
$ws.Range("Q7").Value = 0.1601941747572816
$ws.Range("R7").Value = 0.0970873786407767
$ws.Range("S7").Value = 0.3980582524271845
$ws.Range("B8").Value = 0.06666666666666667
$ws.Range("D8").Value = 0.008888888888888889
$ws.Range("E8").Value = 0.002222222222222222
$ws.Range("F8").Value = 0.08888888888888889
$ws.Range("J8").Value = 0.09555555555555556
$ws.Range("O8").Value = 0.02
$ws.Range("Q8").Value = 0.1955555555555556
$ws.Range("R8").Value = 0.1
$ws.Range("S8").Value = 0.4222222222222222
$ws.Range("B9").Value = 0.07575757575757576
$ws.Range("E9").Value = 0.005050505050505051
$ws.Range("F9").Value = 0.02525252525252525
$ws.Range("J9").Value = 0.06565656565656566
$ws.Range("O9").Value = 0.005050505050505051
$ws.Range("Q9").Value = 0.1919191919191919
$ws.Range("R9").Value = 0.1262626262626263
$ws.Range("S9").Value = 0.5050505050505051
$ws.Range("B10").Value = 0.08536585365853659
$ws.Range("D10").Value = 0.01594746716697936
$ws.Range("E10").Value = 0.0009380863039399625
$ws.Range("F10").Value = 0.0975609756097561
$ws.Range("J10").Value = 0.1060037523452158
$ws.Range("O10").Value = 0.02720450281425891
$ws.Range("Q10").Value = 0.1857410881801126
$ws.Range("R10").Value = 0.09380863039399624
$ws.Range("S10").Value = 0.3874296435272045
$ws.Range("G11").Value = 0.1530612244897959
$ws.Range("J11").Value = 0.07142857142857142
$ws.Range("K11").Value = 0.2074829931972789
$ws.Range("L11").Value = 0.5578231292517006
$ws.Range("S11").Value = 0.01020408163265306
$ws.Range("G12").Value = 0.8117647058823529
$ws.Range("J12").Value = 0.1411764705882353
$ws.Range("K12").Value = 0.005882352941176471
$ws.Range("L12").Value = 0.02941176470588235
$ws.Range("S12").Value = 0.01176470588235294
$ws.Range("G13").Value = 0.7142857142857143
$ws.Range("J13").Value = 0.2619047619047619
$ws.Range("S13").Value = 0.02380952380952381
$ws.Range("J14").Value = 1
$ws.Range("F15").Value = 0.01276595744680851
$ws.Range("H15").Value = 0.1702127659574468
$ws.Range("I15").Value = 0.0851063829787234
$ws.Range("J15").Value = 0.3148936170212766
$ws.Range("K15").Value = 0.07234042553191489
$ws.Range("M15").Value = 0.008510638297872341
$ws.Range("N15").Value = 0.008510638297872341
$ws.Range("O15").Value = 0.05531914893617021
$ws.Range("S15").Value = 0.2723404255319149
$ws.Range("F16").Value = 0.02325581395348837
$ws.Range("H16").Value = 0.2093023255813954
$ws.Range("I16").Value = 0.04651162790697674
$ws.Range("J16").Value = 0.3875968992248062
$ws.Range("K16").Value = 0.06201550387596899
$ws.Range("M16").Value = 0.01550387596899225
$ws.Range("O16").Value = 0.08527131782945736
$ws.Range("S16").Value = 0.1705426356589147
$ws.Range("F17").Value = 0.02233250620347394
$ws.Range("H17").Value = 0.1861042183622829
$ws.Range("I17").Value = 0.109181141439206
$ws.Range("J17").Value = 0.3647642679900744
$ws.Range("K17").Value = 0.09429280397022333
$ws.Range("M17").Value = 0.02481389578163772
$ws.Range("N17").Value = 0.002481389578163772
$ws.Range("O17").Value = 0.06699751861042183
$ws.Range("S17").Value = 0.1290322580645161
$ws.Range("F18").Value = 0.02702702702702703
$ws.Range("H18").Value = 0.2027027027027027
$ws.Range("I18").Value = 0.1261261261261261
$ws.Range("J18").Value = 0.3243243243243243
$ws.Range("K18").Value = 0.1171171171171171
$ws.Range("M18").Value = 0.01801801801801802
$ws.Range("O18").Value = 0.06306306306306306
$ws.Range("S18").Value = 0.1216216216216216
$ws.Range("F19").Value = 0.01973684210526316
$ws.Range("H19").Value = 0.21875
$ws.Range("I19").Value = 0.08141447368421052
$ws.Range("J19").Value = 0.3462171052631579
$ws.Range("K19").Value = 0.1217105263157895
$ws.Range("M19").Value = 0.01973684210526316
$ws.Range("N19").Value = 0.0008223684210526315
$ws.Range("O19").Value = 0.07648026315789473
$ws.Range("S19").Value = 0.1151315789473684
